# Update Sheets via scheduled runner: refresh market-price-derived
# profit figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 855.1818
$ws.Range("I6").Value = 50.875
$ws.Range("K6").Value = 152.625
$ws.Range("M6").Value = -40.625

$ws.Range("H20").Value = 62333.332
$ws.Range("I20").Value = 62333.332
$ws.Range("K20").Value = 62333.332
$ws.Range("M20").Value = -62103.332

$ws.Range("H31").Value = 2122
$ws.Range("I31").Value = 1500
$ws.Range("J31").Value = 2744
$ws.Range("K31").Value = 4500
$ws.Range("L31").Value = 8232
$ws.Range("M31").Value = -4270
$ws.Range("N31").Value = -8692

$ws.Range("H35").Value = 62333.332
$ws.Range("I35").Value = 62333.332
$ws.Range("K35").Value = 62333.332
$ws.Range("M35").Value = -61954.332

$ws.Range("H39").Value = 161.63637
$ws.Range("I39").Value = 70.14286
$ws.Range("J39").Value = 321.75
$ws.Range("K39").Value = 210.42858
$ws.Range("L39").Value = 965.25
$ws.Range("M39").Value = 85.57141999999999
$ws.Range("N39").Value = -1557.25

$ws.Range("H63").Value = 91450
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101248

$ws.Range("H66").Value = 91450
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306240

$ws.Range("H86").Value = 898.75
$ws.Range("I86").Value = 898.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 898.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 224.25
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 898.75
$ws.Range("I89").Value = 898.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 4493.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 1122.25
$ws.Range("N89").Value = ""

$ws.Range("H113").Value = 10107.8
$ws.Range("I113").Value = 8884.5
$ws.Range("K113").Value = 8884.5
$ws.Range("M113").Value = -5630.5

$ws.Range("H132").Value = 4499.9287
$ws.Range("I132").Value = 4750.25
$ws.Range("K132").Value = 14250.75
$ws.Range("M132").Value = -11720.75

$ws.Range("H138").Value = 2194.8125
$ws.Range("J138").Value = 2882.111
$ws.Range("L138").Value = 8646.332999999999
$ws.Range("N138").Value = -18926.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1596.2653
$ws.Range("I74").Value = 1532.9231
$ws.Range("J74").Value = 1843.3
$ws.Range("K74").Value = 1532.9231
$ws.Range("L74").Value = 1843.3
$ws.Range("M74").Value = -658.9231
$ws.Range("N74").Value = -3591.3

$ws.Range("H77").Value = 1596.2653
$ws.Range("I77").Value = 1532.9231
$ws.Range("J77").Value = 1843.3
$ws.Range("K77").Value = 7664.6155
$ws.Range("L77").Value = 9216.5
$ws.Range("M77").Value = -3296.6155
$ws.Range("N77").Value = -17952.5

$ws.Range("H110").Value = 1464.2
$ws.Range("I110").Value = 1530.6666
$ws.Range("J110").Value = 1198.3334
$ws.Range("K110").Value = 1530.6666
$ws.Range("L110").Value = 1198.3334
$ws.Range("M110").Value = 514.3334
$ws.Range("N110").Value = -5288.3334

$ws.Range("H122").Value = 4753.0566
$ws.Range("I122").Value = 5536.364
$ws.Range("K122").Value = 16609.092
$ws.Range("M122").Value = -14159.092

$ws.Range("H132").Value = 7687.164
$ws.Range("I132").Value = 7717.912
$ws.Range("J132").Value = 7249
$ws.Range("K132").Value = 23153.736
$ws.Range("L132").Value = 21747
$ws.Range("M132").Value = -20623.736
$ws.Range("N132").Value = -26807

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 33637
$ws.Range("I26").Value = 33637
$ws.Range("K26").Value = 33637
$ws.Range("M26").Value = -33345

$ws.Range("H86").Value = 2734.353
$ws.Range("I86").Value = 1718.8889
$ws.Range("J86").Value = 3876.75
$ws.Range("K86").Value = 1718.8889
$ws.Range("L86").Value = 3876.75
$ws.Range("M86").Value = -595.8888999999999
$ws.Range("N86").Value = -6122.75

$ws.Range("H89").Value = 2734.353
$ws.Range("I89").Value = 1718.8889
$ws.Range("J89").Value = 3876.75
$ws.Range("K89").Value = 8594.4445
$ws.Range("L89").Value = 19383.75
$ws.Range("M89").Value = -2978.4445
$ws.Range("N89").Value = -30615.75

$ws.Range("H96").Value = 9843.833000000001
$ws.Range("I96").Value = 9843.833000000001
$ws.Range("K96").Value = 9843.833000000001
$ws.Range("M96").Value = -7097.833000000001

$ws.Range("H105").Value = 2873.923
$ws.Range("I105").Value = 3881.5334
$ws.Range("K105").Value = 3881.5334
$ws.Range("M105").Value = -2134.5334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4483.1665
$ws.Range("I58").Value = 4966.6665
$ws.Range("J58").Value = 3999.6667
$ws.Range("K58").Value = 4966.6665
$ws.Range("L58").Value = 3999.6667
$ws.Range("M58").Value = -4763.6665
$ws.Range("N58").Value = -4405.6667

$ws.Range("H107").Value = 1111.4445
$ws.Range("I107").Value = 690.875
$ws.Range("K107").Value = 690.875
$ws.Range("M107").Value = 1229.125

$ws.Range("H136").Value = 4483.1665
$ws.Range("I136").Value = 4966.6665
$ws.Range("J136").Value = 3999.6667
$ws.Range("K136").Value = 14899.9995
$ws.Range("L136").Value = 11999.0001
$ws.Range("M136").Value = -12349.9995
$ws.Range("N136").Value = -17099.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1622.92
$ws.Range("J68").Value = 2331.6667
$ws.Range("L68").Value = 6995.000100000001
$ws.Range("N68").Value = -8617.000100000001

$ws.Range("H71").Value = 1622.92
$ws.Range("J71").Value = 2331.6667
$ws.Range("L71").Value = 20985.0003
$ws.Range("N71").Value = -29097.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3745.818
$ws.Range("I113").Value = 3295.6
$ws.Range("J113").Value = 4121
$ws.Range("K113").Value = 3295.6
$ws.Range("L113").Value = 4121
$ws.Range("M113").Value = -1125.6
$ws.Range("N113").Value = -8461

$ws.Range("H132").Value = 5521.1816
$ws.Range("I132").Value = 5476.1113
$ws.Range("K132").Value = 16428.3339
$ws.Range("M132").Value = -13898.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62505690
$ws.Range("I7").Value = 5916.5
$ws.Range("K7").Value = 5916.5
$ws.Range("M7").Value = -5804.5

$ws.Range("H46").Value = 855
$ws.Range("J46").Value = 882.4
$ws.Range("N46").Value = -1258.4

$ws.Range("H126").Value = 62505690
$ws.Range("I126").Value = 5916.5
$ws.Range("K126").Value = 17749.5
$ws.Range("M126").Value = -15279.5
